# Main.xlsx / "Rules" sheet: rename rule R40 (row 11) to "1".
#
# B11 must stay a TEXT cell (it is the "Rule" name column), so we can't
# just assign Value = "1" directly -- Excel's smart-typing would store
# that as the number 1 instead of the string "1". To force text without
# disturbing the cell's existing style/format (no quote-prefix style,
# no NumberFormat="@" residue), stage the text in a throwaway cell that
# has been explicitly formatted as Text, copy it, and paste only the
# value into B11 -- PasteSpecial(xlPasteValues) leaves B11's own
# formatting (border/fill/font/numFmt) completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()
$target.PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()
